# Apply weekly update to the "Hortaliza, Femacal de La Calera - Ciboulette" sheet.
# The edit inserts a new weekly record and shifts the existing date/price series
# for rows 102-291 down by one (each row takes on the values previously held by
# the row above it), while row 101 receives the newest record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 101
$lastOldRow = 290
$newLastRow = 291

# Columns that vary per-row and must be shifted: D (date), J (price), K, L, M (range
# prices) and P (per-unit price). All other columns (A, B, C, E, F, G, H, I, N, O, Q, R)
# are constant for every data row in this sheet, so they do not need to be touched,
# except for populating the brand-new row 291 which does not exist yet.

# 1) Capture the existing values for the columns that will shift, for rows
#    firstDataRow..lastOldRow, before making any changes.
$colD = @{}
$colJ = @{}
$colK = @{}
$colL = @{}
$colM = @{}
$colP = @{}

for ($i = $firstDataRow; $i -le $lastOldRow; $i++) {
    $colD[$i] = $ws.Cells.Item($i, 4).Value()
    $colJ[$i] = $ws.Cells.Item($i, 10).Value()
    $colK[$i] = $ws.Cells.Item($i, 11).Value()
    $colL[$i] = $ws.Cells.Item($i, 12).Value()
    $colM[$i] = $ws.Cells.Item($i, 13).Value()
    $colP[$i] = $ws.Cells.Item($i, 16).Value()
}

# 2) Populate the brand-new row 291 with the constant columns (same as every
#    other data row on this sheet) before anything else touches it.
$ws.Cells.Item($newLastRow, 1).Value = 3
$ws.Cells.Item($newLastRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newLastRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newLastRow, 5).Value = 5
$ws.Cells.Item($newLastRow, 6).Value = 100112039
$ws.Cells.Item($newLastRow, 7).Value = "Ciboulette"
$ws.Cells.Item($newLastRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newLastRow, 9).Value = "Primera"
$ws.Cells.Item($newLastRow, 14).Value = "`$/docena de atados"
$ws.Cells.Item($newLastRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newLastRow, 17).Value = 3
$ws.Cells.Item($newLastRow, 18).Value = "Hortaliza"

# 3) Shift rows 102..291 so that row i gets what row (i-1) used to hold.
for ($i = $newLastRow; $i -ge ($firstDataRow + 1); $i--) {
    $src = $i - 1
    $ws.Cells.Item($i, 4).Value = $colD[$src]
    $ws.Cells.Item($i, 10).Value = $colJ[$src]
    $ws.Cells.Item($i, 11).Value = $colK[$src]
    $ws.Cells.Item($i, 12).Value = $colL[$src]
    $ws.Cells.Item($i, 13).Value = $colM[$src]
    $ws.Cells.Item($i, 16).Value = $colP[$src]
}

# 4) Fill in the brand-new newest record at row 101.
$ws.Cells.Item($firstDataRow, 4).Value = 44645
$ws.Cells.Item($firstDataRow, 10).Value = 120

# 5) Make sure the date column in the new row uses the same date/time number
#    format as the rest of column D.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastOldRow, 4).NumberFormat()

Write-Output "done"
